# Update "想去人数" (want-to-go count) figures in the F column of the
# 展览 (Exhibitions) and 全部类型 (All types) sheets, which mirror the
# same underlying data.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 7288
    4  = 5514
    8  = 43
    9  = 106
    12 = 204
    13 = 45
    15 = 293
    19 = 41
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
